$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "nest" field: rename the existing 2Dlist header and extend its type
#     definition on row 5 with a second nested list `[[bool]]`.
$ws.Range("V6").Value = "nest"

$ws.Range("AB5").Value = "["
$ws.Range("AC5").Value = "["
$ws.Range("AD5").Value = "bool"
$ws.Range("AE5").Value = "]"
$ws.Range("AF5").Value = "]"
$ws.Range("AG5").Value = "]"

# Sample data row (row 7): boolean value "Y" under the new bool column.
$ws.Range("AD7").Value = "Y"

# Update selection to reflect where the editor ended up.
$ws.Range("AG10").Select()
